$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 33 currently carries the "last row" styling (the fill-flagged variant
# of the date/number/text formats). A new row is being appended below it, so
# row 33 becomes a normal interior row and the new row 34 becomes the last
# row, inheriting the styling row 33 used to have.

# 1) Propagate row 33's current formatting down into the new row 34 first,
#    before row 33's own formatting changes.
$ws.Range("A33:C33").Copy()
$ws.Range("A34:C34").PasteSpecial(-4122)

# 2) Restore row 33 to the regular (non-last-row) formatting by copying the
#    format down from row 32, which already uses the normal styling.
$ws.Range("A32:C32").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Populate the new row's values.
$ws.Range("A34").Value = 45728
$ws.Range("B34").Value = 6
$ws.Range("C34").Value = "Finalized Network Security Scanning Tool : Nmap, Web App Security Scanning Tool : OWASP ZAP. Another verification for OWASP ZAP is required. If this can not be verified, I need select second best tool."

# 4) Match the row height used by the other wrapped, multi-line rows.
$ws.Rows.Item(34).RowHeight = 31.5
